$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("MANDATORY"); this shifts the old
# C..J (MANDATORY..ACCEPTED) one column right to D..K, matching the diff.
$ws.Range("C1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "DEFAULT"

# New column's width (~8.66 characters, matching the author's resize).
$ws.Range("C1").ColumnWidth = 7.833333333333333

# Default values for the new DEFAULT column (all FALSE booleans).
$ws.Range("C2:C5").Value = $false

# Match cell formatting of the neighbouring DESCRIPTION column (black font,
# same style bucket used elsewhere for data-row cells) for the new column's
# data cells - mirrors the s="3" style seen for B2:B3/F3:F5 in the sheet.
$ws.Range("C2:C5").Font.Color = 0

# Row 13 gains an (empty) formatted cell under the new column too.
$ws.Range("C13").Font.Color = 0

# Update the selected cell shown when the workbook is reopened.
$ws.Range("F10").Select()

# Page setup: paper size + orientation now explicit (A4, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
